$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B descriptions are being re-ordered / replaced for rows 2-26 (FACULTY1..FACULTY25
# test files got new / re-shuffled descriptions, two obsolete descriptions were dropped from
# the shared string table and seven new "Line one ..." descriptions were appended).
# Row order below matches the order the new strings were appended to the shared-string table.
$rowOrder = @(3, 2, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26)

$descriptions = @{
    2  = "Line 1 input rejected. Decimal numbers are not acceptable."
    3  = "Line one input rejected minimum teaching hours is 1."
    4  = "No space after comma that follows the last name of the faculty member. The input is rejected."
    5  = "Input rejected. Only white space can be used as space seperator."
    6  = "Input rejected because email address is not provided."
    7  = "Input rejected because hours teaching is not provided."
    8  = "Input rejected because years of service is not provided."
    9  = "Input rejected because first and last name is missing from input."
    10 = "Input rejected. First name of faculty is not all upper-case."
    11 = "Input rejected. Years of service exceeds the maxiumum of  60."
    12 = "Input rejected. Empty file."
    13 = "Line one input rejected because email extension is partially lowercase"
    14 = "Input rejected. Years of service cannot be a decimal number."
    15 = "Line one input rejected. Email length exceeds the maximum of 10."
    16 = "line one input is rejected. Email length is less than the minimum 3."
    17 = "Input should be accepted. The years of service value is 60, and the maximum is also 60."
    18 = "Input should be accepted. The minimum hours teaching is three."
    19 = "Input should be accepted. Medium to large file input."
    20 = "Input should be accepted. Large amount of input."
    21 = "All file input is not accepted. Line 1 minimum hours is less than three. Line two email extension is partially lower-case. Line three email address is lower-case. Line four only part of email address is lower-case."
    22 = "File contents not accepted. Line one last name is lower-case. Line two first name is lower-case, and line three is missing a space between last and first name."
    23 = "Input is not accepted because the `"@UNA.EDU`" is missing from email. The name has the number `"12`" preceeding it."
    24 = "The input is not accepted. On the first line of input, a backslash follows the comma after the last name. The email address is missing the `"EDU`" extension. On the second line the `"@`" symbol is missing from the email."
    25 = "The input is not accepted. The first two lines are duplicates. The third line the years of service is less than zero, and the hourd teaching is less than 3."
    26 = "The first line in the input file is missing the years of service. The second line in the input file has the same email address for a different faulty member."
}

foreach ($row in $rowOrder) {
    $ws.Cells.Item($row, 2).Value = $descriptions[$row]
}

# Row 22's custom height grew to accommodate its new (longer) description.
$ws.Rows.Item(22).RowHeight = 56.25

# Reflect the author's final selection/scroll position on Sheet1.
$ws.Range("B23").Select()

$wb.Save()
